$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix bug when exporting an array of native types:
# struct B's field "d:uint8" should be "d:[]uint8"
$ws.Range("E5").Value = "d:[]uint8"

# Update selection (view state) as recorded in the saved file
$ws.Range("G11").Select()
